$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use row 66 (an existing fully-styled data row) as the style template so the
# new rows' cells pick up the same styles (s="2" for the date column, s="1"
# for everything else) instead of getting a brand-new style index.
$templateRow = 66

function Set-CellWithTemplateStyle {
    param(
        [string]$Col,
        [int]$Row,
        $Value
    )
    $ws.Range($Col + $templateRow).Copy() | Out-Null
    $ws.Range($Col + $Row).PasteSpecial(-4122) | Out-Null
    $ws.Range($Col + $Row).Value = $Value
}

# ---- Row 67 ----
$row = 67
Set-CellWithTemplateStyle "A" $row 43550.947415567134
Set-CellWithTemplateStyle "B" $row "Depuis + de 8 ans"
Set-CellWithTemplateStyle "C" $row "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Livres, Vidéos Youtube de vulgarisation, Articles de presse, Publications scientifiques"
Set-CellWithTemplateStyle "D" $row 5.0
Set-CellWithTemplateStyle "E" $row 9.0
Set-CellWithTemplateStyle "F" $row "A déjà commencé"
Set-CellWithTemplateStyle "G" $row "Assez lent (de l'ordre de 50 ans et plus)"
Set-CellWithTemplateStyle "H" $row "Angoisse"
Set-CellWithTemplateStyle "I" $row "On a un fort potentiel d'action à l'échelle individuelle, Je suis prêt à baisser mon niveau de vie si cette baisse s'opère pour les autres également, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique, Pour agir efficacement, il faut hiérarchiser les problèmes (perte de biodiversité, dérèglement climatique etc.)"
Set-CellWithTemplateStyle "J" $row "Une réaction de personnes ne voulant pas toucher à leur mode de vie"
Set-CellWithTemplateStyle "K" $row 8.0
Set-CellWithTemplateStyle "L" $row 6.0
Set-CellWithTemplateStyle "M" $row 7.0
Set-CellWithTemplateStyle "N" $row 5.0
Set-CellWithTemplateStyle "O" $row 7.0
Set-CellWithTemplateStyle "P" $row "Je partage des liens sur les réseaux sociaux, Je fais profil bas. Trop en parler, c'est devenir prêcheur, et donc desservir la cause., Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
Set-CellWithTemplateStyle "Q" $row 3.0
Set-CellWithTemplateStyle "S" $row "Un homme"
Set-CellWithTemplateStyle "T" $row 27.0
Set-CellWithTemplateStyle "U" $row "En ville dans une grande agglomération"
Set-CellWithTemplateStyle "V" $row "Doctorat ou équivalent"
Set-CellWithTemplateStyle "W" $row "Génie électrique"
Set-CellWithTemplateStyle "X" $row "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
Set-CellWithTemplateStyle "Y" $row "Flexitarien"
Set-CellWithTemplateStyle "Z" $row "Très à gauche (France insoumise ou plus à gauche en France)"
Set-CellWithTemplateStyle "AB" $row "Coach, Ancien Académicien"

# ---- Row 68 ----
$row = 68
Set-CellWithTemplateStyle "A" $row 43551.650507453705
Set-CellWithTemplateStyle "B" $row "Depuis + de 8 ans"
Set-CellWithTemplateStyle "C" $row "Articles de vulgarisation & blogs, Vidéos Youtube de vulgarisation, Articles de presse, Publications scientifiques"
Set-CellWithTemplateStyle "D" $row 4.0
Set-CellWithTemplateStyle "E" $row 8.0
Set-CellWithTemplateStyle "F" $row "Va commencer dans les 15 à 20 ans qui viennent"
Set-CellWithTemplateStyle "G" $row "Un peu plus lent (de l'ordre de 20 à 30 ans)"
Set-CellWithTemplateStyle "H" $row "Angoisse"
Set-CellWithTemplateStyle "I" $row "On a un fort potentiel d'action à l'échelle individuelle, Je comprend qu'il y ait des personnes climatosceptiques au sein de la population, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique"
Set-CellWithTemplateStyle "J" $row "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Problème d'éducation et/ou d'information"
Set-CellWithTemplateStyle "K" $row 7.0
Set-CellWithTemplateStyle "L" $row 6.0
Set-CellWithTemplateStyle "M" $row 7.0
Set-CellWithTemplateStyle "N" $row 2.0
Set-CellWithTemplateStyle "O" $row 5.0
Set-CellWithTemplateStyle "P" $row "Je partage des liens sur les réseaux sociaux, Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi), Je montre l'exemple, je suscite l'étonnement"
Set-CellWithTemplateStyle "Q" $row 3.0
Set-CellWithTemplateStyle "R" $row "Une religion c'est entre autre fait pour expliquer la mort, la vie après, la réincarnation etc. Avec l'écologie on parle de tout sauf d'une vie meilleure après la mort."
Set-CellWithTemplateStyle "S" $row "Une femme"
Set-CellWithTemplateStyle "T" $row 24.0
Set-CellWithTemplateStyle "U" $row "En ville dans une grande agglomération"
Set-CellWithTemplateStyle "V" $row "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
Set-CellWithTemplateStyle "W" $row "Mathématiques / Informatique"
Set-CellWithTemplateStyle "X" $row "Très frugal (flexitarien ou végétalien, AMAP, déplacement doux)"
Set-CellWithTemplateStyle "Y" $row "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
Set-CellWithTemplateStyle "Z" $row "Parti à préoccupation environnementale (Europe Écologie les Verts en France)"
Set-CellWithTemplateStyle "AB" $row "Académicien"
